# admin panel system change
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be created in the same order they first appear
# in the final workbook (matches sharedStrings.xml append order).

# Row 5: add C5 = new string, row height already 75 (unchanged)
$ws.Range("C5").Value = "وسط چین بودن عکس های بلاگ"

# Row 6: add C6 = new string, row height grows from 30 to 75
$ws.Range("C6").Value = "مطالب پیشنهادی وبلاگ بیش از 6 است"
$ws.Rows(6).RowHeight = 75

# Row 4: C4 changes from "courses" to new "admin panel get global"
$ws.Range("C4").Value = "admin panel get global"

# Row 21: add B21 = new string
$ws.Range("B21").Value = "یادآور ثبت محتوا"

# Row 22: add B22 = new string
$ws.Range("B22").Value = "زمانبندی ارسال محتوا"

# Row 17: add B17 = "courses" (moved here from old C4, reuses existing string)
$ws.Range("B17").Value = "courses"

# Update view: scroll so row 19 is at the top of the window, select B23
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
